$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update localization status text: "Ready for handoff" -> "In Translation"
#    on every sheet that shows a Status column.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Replace("Ready for handoff", "In Translation")

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Replace("Ready for handoff", "In Translation")

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Replace("Ready for handoff", "In Translation")

# ---------------------------------------------------------------------------
# 2. Narrow the "Status" columns (previously auto-sized for the longer
#    "Ready for handoff" text) now that the shorter "In Translation" text
#    fits in a narrower column.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
